$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.169.26"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").Value = "1.589.26"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'211.72"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").Value = "'0.501"
$ws.Range("E6").Value = "  -1.02%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").Value = "'0.0604"
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").Value = "'18.97"
$ws.Range("E10").Value = "  -2.21%  "
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "1.813.97"
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").Value = "1.581.38"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("E14").Value = "  -1.53%  "
$ws.Range("D15").Value = "'0.509"
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "'63.55"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "26.179.55"
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "0.0₃0723"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("D19").Value = "'214.50"
$ws.Range("E19").Value = "  +1.91%  "
$ws.Range("D20").Value = "'7.34"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "'4.23"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "'9.00"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("D25").Value = "'144.79"
$ws.Range("E25").Value = "  +0.13%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "'6.95"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("E28").Value = "  -1.11%  "
$ws.Range("D29").Value = "'15.05"
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("E30").Value = "  -2.49%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  -1.66%  "
$ws.Range("D33").Value = "1.417.24"
$ws.Range("E33").Value = "  +8.01%  "
$ws.Range("E34").Value = "  -1.67%  "
$ws.Range("D35").Value = "'2.42"
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").Value = "'0.586"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").Value = "'0.823"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("E40").Value = "  +4.87%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'0.966"
$ws.Range("E42").Value = "  -11.18%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "'0.764"
$ws.Range("E44").Value = "  -0.45%  "
$ws.Range("D45").Value = "1.725.45"
$ws.Range("D46").Value = "'60.93"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("D47").Value = "'86.91"
$ws.Range("E47").Value = "  -0.70%  "
$ws.Range("E48").Value = "  -0.82%  "
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("D50").Value = "'0.0960"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("E51").Value = "  -0.16%  "
